# "feedback form and other component updated"
#
# - rename Sheet2 -> Student, make it the active tab (was FacultyTbl active)
# - build out the Student sheet: header row + one student record, with a
#   mailto hyperlink on the e-mail cell and matching column widths
# - page orientation on the Student sheet

$wb = $excel.ActiveWorkbook
$student = $wb.Worksheets.Item("Sheet2")

$student.Name = "Student"

# --- headers (row 1) + data (row 2) ----------------------------------
# Written in the same order the original author typed them in Excel, so
# the shared-string table comes out in the same append order.
$student.Range("K1").Value = "Address"
$student.Range("K2").Value = "sdjfhbhjb"
$student.Range("C1").Value = "FathersName"
$student.Range("D1").Value = "MothersName"
$student.Range("G1").Value = "FathersContactNo"
$student.Range("H1").Value = "MothersContactNo"
$student.Range("I1").Value = "ParentsLoginPassword"
$student.Range("A1").Value = "RollNo"
$student.Range("A2").Value = "16EGICS039"
$student.Range("B2").Value = "Himanshu Panchal"
$student.Range("C2").Value = "Mr. Rahul Panchal"
$student.Range("D2").Value = "Mrs. Jaya Panchal"
$student.Range("E2").Value = "panchalhimanshu@gmail.com"

$student.Range("B1").Value = "Name"
$student.Range("E1").Value = "EmailID"
$student.Range("F1").Value = "ContactNo"
$student.Range("J1").Value = "Password"
# A1 (RollNo) stays un-bold; only B1:K1 are bold headers
$student.Range("B1:K1").Font.Bold = $true

$student.Range("F2").Value = 7410258963
$student.Range("G2").Value = 9874563210
$student.Range("H2").Value = 8520369147
$student.Range("I2").Value = 111111
$student.Range("J2").Value = 222222

# e-mail hyperlink (row 2) + a matching blank styled cell on row 3
$student.Hyperlinks.Add($student.Range("E2"), "mailto:panchalhimanshu@gmail.com")
$student.Range("E2").Style = "Hyperlink"
$student.Range("E3").Style = "Hyperlink"

# --- column widths -------------------------------------------------------
$student.Columns.Item(1).ColumnWidth = 17.736979166666668
$student.Columns.Item(2).ColumnWidth = 19.166666666666668
$student.Columns.Item(3).ColumnWidth = 23.451822916666668
$student.Columns.Item(4).ColumnWidth = 22.592447916666668
$student.Columns.Item(5).ColumnWidth = 31.451822916666668
$student.Columns.Item(6).ColumnWidth = 18.307291666666668
$student.Columns.Item(7).ColumnWidth = 18.307291666666668
$student.Columns.Item(8).ColumnWidth = 18.307291666666668
$student.Columns.Item(9).ColumnWidth = 18.307291666666668
$student.Columns.Item(10).ColumnWidth = 14.736979166666666
$student.Columns.Item(11).ColumnWidth = 42.736979166666664

# --- page setup ------------------------------------------------------------
$student.PageSetup.Orientation = 1

# --- view state: Student tab becomes the active / selected sheet ----------
$student.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$student.Range("L2").Select()
